# Add two new Darwin Core fields to the occurrences template header row:
#   AB1 -> occurrenceID
#   AC1 -> comments
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB1").Value = "occurrenceID"
$ws.Range("AC1").Value = "comments"

# Move the selection to the newly added cell below the new "comments" column,
# and scroll the view so column Q is the left-most visible column.
$ws.Range("AC2").Select() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 17
$aw.ScrollRow = 1
